# Lecture partielle de l'EDT M1 MIAGE.
# Shift all the date entries in column A forward by 1096 days (so the
# schedule is re-used ~3 years later) and update the French weekday
# labels in column B so that they keep matching the (now shifted) dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that hold a date in column A together with the matching weekday
# name in column B. Each date is pushed forward by exactly 1096 days.
$rows = 2, 6, 10, 14, 17, 20, 23, 26, 30, 34, 37, 41, 45, 51

foreach ($r in $rows) {
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $aCell.Value2() + 1096
}

# French weekday names for the updated dates.
$weekdays = @{
    2  = "jeudi"
    6  = "jeudi"
    10 = "lundi"
    14 = "mercredi"
    17 = "jeudi"
    20 = "vendredi"
    23 = "mardi"
    26 = "mercredi"
    30 = "mardi"
    34 = "lundi"
    37 = "mercredi"
    41 = "jeudi"
    45 = "lundi"
    51 = "lundi"
}

foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $weekdays[$r]
}
